$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-11-20 Wednesday" "2024-11-21 Thursday"

Replace-Text "65×63=" "35×72="
Replace-Text "93×78=" "84×48="
Replace-Text "60×27=" "71×37="
Replace-Text "62×18=" "11×42="
Replace-Text "64×11=" "68×50="
Replace-Text "54×18=" "33×18="
Replace-Text "79×32=" "12×30="
Replace-Text "22×29=" "89×50="
Replace-Text "67×33=" "15×68="
Replace-Text "27×65=" "99×62="
Replace-Text "70×24=" "66×39="
Replace-Text "32×81=" "60×24="
Replace-Text "70×90=" "28×76="
Replace-Text "38×65=" "44×84="
Replace-Text "34×99=" "90×23="
Replace-Text "11×15=" "19×35="
Replace-Text "39×76=" "29×48="
Replace-Text "97×53=" "17×38="
Replace-Text "15×27=" "45×16="
Replace-Text "69×67=" "28×26="
Replace-Text "43×12=" "74×20="
Replace-Text "33×44=" "86×80="
Replace-Text "40×71=" "58×87="
Replace-Text "64×55=" "36×11="
Replace-Text "72×26=" "14×48="
